$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.948.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.673.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.86%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.910.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.672.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.968.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("E19").Value = "  +3.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0731"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("E22").Value = "  -1.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.22%  "

$ws.Range("E24").Value = "  -2.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("E31").Value = "  -0.67%  "

$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.483.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.42%  "

$ws.Range("E35").Value = "  +3.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("E37").Value = "  +1.51%  "

$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.894"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("E40").Value = "  -3.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.97%  "

$ws.Range("E43").Value = "  +1.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.816.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.777"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("E49").Value = "  +1.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0510"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("E51").Value = "  +0.21%  "
